$d = $word.ActiveDocument

# The "Ancillary Structures" page used to embed a (placeholder) picture in
# the FirstParagraph right under the heading. It is replaced by a plain
# hyperlink whose display text is the picture's original URL on
# ura.gov.sg, styled with the document's "Hyperlink" character style.
$imageUrl = "https://ura.gov.sg/-/media/Corporate/Guidelines/Development-control/Others/HMC03_Ancillary_Structures.jpg?h=100%25&w=100%25"

$shape = $d.InlineShapes.Item(1)
$shapeRange = $shape.Range
$shapeRange.Text = $imageUrl
$d.Hyperlinks.Add($shapeRange, $imageUrl)
